# Update countries & provincias Spain
# - Swap the display order of three adjacent country rows (Bielorrusia/Ecuador,
#   Croacia/Consejo Danes para los Refugiados, Montserrat/Seychelles)
# - Refresh the statistics (B:H) for the affected rows with the newer figures
# - Bump the "Datos actualizados" timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp -------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 19:35"

# --- Row 4: Estados Unidos (values refreshed only) ----------------------
$ws.Cells.Item(4, 2).Value = 1694562
$ws.Cells.Item(4, 3).Value = 8126
$ws.Cells.Item(4, 4).Value = 456397
$ws.Cells.Item(4, 5).Value = 1138701
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 164
$ws.Cells.Item(4, 8).Value = 99464

# --- Row 12: Turquia (values refreshed only) -----------------------------
$ws.Cells.Item(12, 2).Value = 157814
$ws.Cells.Item(12, 3).Value = 987
$ws.Cells.Item(12, 4).Value = 120015
$ws.Cells.Item(12, 5).Value = 33430
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 29
$ws.Cells.Item(12, 8).Value = 4369

# --- Row 13: India (values refreshed only) -------------------------------
$ws.Cells.Item(13, 2).Value = 144741
$ws.Cells.Item(13, 3).Value = 6205
$ws.Cells.Item(13, 4).Value = 60498
$ws.Cells.Item(13, 5).Value = 80081
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 138
$ws.Cells.Item(13, 8).Value = 4162

# --- Rows 25/26: swap Bielorrusia/Ecuador order + refresh values --------
$ws.Cells.Item(25, 1).Value = "Ecuador"
$ws.Cells.Item(25, 2).Value = 37355
$ws.Cells.Item(25, 3).Value = 599
$ws.Cells.Item(25, 4).Value = 3560
$ws.Cells.Item(25, 5).Value = 30592
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 95
$ws.Cells.Item(25, 8).Value = 3203

$ws.Cells.Item(26, 1).Value = "Bielorrusia"
$ws.Cells.Item(26, 2).Value = 37144
$ws.Cells.Item(26, 3).Value = 946
$ws.Cells.Item(26, 4).Value = 14449
$ws.Cells.Item(26, 5).Value = 22491
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 5
$ws.Cells.Item(26, 8).Value = 204

# --- Row 42: Israel (values refreshed only) ------------------------------
$ws.Cells.Item(42, 2).Value = 16734
$ws.Cells.Item(42, 3).Value = 17
$ws.Cells.Item(42, 4).Value = 14307
$ws.Cells.Item(42, 5).Value = 2146
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 2
$ws.Cells.Item(42, 8).Value = 281

# --- Rows 86/87: swap Croacia/Consejo Danes order + refresh values ------
$ws.Cells.Item(86, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(86, 2).Value = 2297
$ws.Cells.Item(86, 3).Value = 156
$ws.Cells.Item(86, 4).Value = 337
$ws.Cells.Item(86, 5).Value = 1893
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 4
$ws.Cells.Item(86, 8).Value = 67

$ws.Cells.Item(87, 1).Value = "Croacia"
$ws.Cells.Item(87, 2).Value = 2244
$ws.Cells.Item(87, 3).Value = 0
$ws.Cells.Item(87, 4).Value = 2035
$ws.Cells.Item(87, 5).Value = 109
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 1
$ws.Cells.Item(87, 8).Value = 100

# --- Rows 210/211: swap Montserrat/Seychelles order + refresh values ----
$ws.Cells.Item(210, 1).Value = "Seychelles"
$ws.Cells.Item(210, 2).Value = 11
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 2).Value = 11
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 1
